# edit.ps1
# Implements:
#  - rename defined name Master_Gate_Status -> ComplianceGate
#  - rename defined name Compliance_Adjusted_Value -> AdjustedValue
#  - propagate the rename through every formula / label that uses those names
#  - update a handful of gate-logic cells on 11_Calc_Compliance (NAFC-based
#    toxicity gate, DO gate operator change)
#  - change the Net Value formula on 14_Calc_Sim (H10) to use Gross_Value
#    instead of the (renamed) AdjustedValue
#  - delete row 19 ("Compliance Gate Status") from 20_Dashboard, shifting
#    everything below it up by one row (tables / merged cells / dimension
#    follow automatically)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename the defined names FIRST (so that cached formula values that
#    are recalculated while we touch cells below resolve correctly),
#    then sweep every worksheet replacing the literal text everywhere
#    it appears (formulas, inline-string labels, notes, etc.)
# ---------------------------------------------------------------------

$renames = @{
    "Master_Gate_Status"        = "ComplianceGate"
    "Compliance_Adjusted_Value" = "AdjustedValue"
}

foreach ($n in $wb.Names) {
    if ($renames.ContainsKey($n.Name)) {
        $n.Name = $renames[$n.Name]
    }
}

# ---------------------------------------------------------------------
# 2. Sheet-specific manual edits that are NOT simple renames
# ---------------------------------------------------------------------

# 11_Calc_Compliance : Toxicity gate now checks NAFC vs target (was a
# stochastic-efficiency baseline check), and the DO gate becomes >=.
$wsCompliance = $wb.Worksheets.Item("11_Calc_Compliance")
$wsCompliance.Range("B8").Value = "NAFC <= Target"
$wsCompliance.Range("C8").Value = "below"
$wsCompliance.Range("D8").Formula = "=IF(Final_NAFC<=Target_NAFC,1,0)"
$wsCompliance.Range("F8").Value = "NAFC toxicity target"

$wsCompliance.Range("B11").Value = ">= 6.5 mg/L"
$wsCompliance.Range("D11").Formula = "=IF(Env_DO_Typical_Jul>=6.5,1,0)"

# 14_Calc_Sim : Net value KPI now compares Gross_Value (not the gated
# value) against Testing_Cost.
$wsSim = $wb.Worksheets.Item("14_Calc_Sim")
$wsSim.Range("H10").Formula = "=Gross_Value-Testing_Cost"

# ---------------------------------------------------------------------
# 3. Propagate the two renames textually across every worksheet (this
#    fixes up formulas that reference the defined names as well as any
#    inline-string cell labels / notes that spell the names out).
# ---------------------------------------------------------------------

foreach ($ws in $wb.Worksheets) {
    foreach ($old in $renames.Keys) {
        $new = $renames[$old]
        $ws.Cells.Replace($old, $new) | Out-Null
    }
}

# ---------------------------------------------------------------------
# 4. 20_Dashboard : remove the now-redundant "Compliance Gate Status"
#    row; everything below shifts up automatically (tables, merged
#    cells, sheet dimension all follow the row delete).
# ---------------------------------------------------------------------

$wsDash = $wb.Worksheets.Item("20_Dashboard")
$wsDash.Rows.Item(19).Delete()
